$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match header style/formatting used by the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2..71 (columns I = I0, J = IF)
$iValues = @(9,8,8,9,9,9,8,9,9,8,7,9,8,6,8,8,7,8,9,8,9,8,9,9,8,9,9,9,9,9,8,8,9,9,9,9,9,9,9,9,8,9,9,8,9,9,8,7,9,8,7,7,8,9,9,7,7,3,7,7,8,8,6,7,9,5,5,5,3,7)
$jValues = @(9,8,8,9,9,9,8,9,9,8,8,9,8,6,8,8,8,8,9,8,9,8,9,9,8,9,9,9,9,9,8,9,9,9,9,9,9,9,9,9,8,9,9,8,9,9,8,7,9,8,7,7,8,9,9,7,7,4,7,7,8,8,6,7,9,6,5,5,3,7)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
